$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set Debug mode value in OSS2Params (G1) from 0.8 to 0.2
$ws.Range("G1").Value2 = 0.2

# Update the mirrored static bound columns (N/O) to match recalculated F/G columns
for ($r = 2; $r -le 20; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $f = $a - 0.2 * [Math]::Abs($a)
    $g = $a + 0.2 * [Math]::Abs($a)
    $ws.Cells.Item($r, 14).Value2 = $f
    $ws.Cells.Item($r, 15).Value2 = $g
}

# Update the selected cell/range shown in the sheet view
$ws.Range("T13").Select()
